$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a zero-padded / numeric-looking account number into a cell as
# TEXT, without Excel's COM auto-conversion turning it into a real number
# (which would drop the leading zeros) and without leaving a stray
# "quote prefix" style on the cell (which the source file's cells don't have).
# Trick: compute the padded text with TEXT() in an unused staging cell, then
# .Copy() that computed value into the destination (copies as a clean text
# value), then clear the staging cell.
# ---------------------------------------------------------------------------
function Set-TextAccount($row, $col, $digits) {
    $stage = $ws.Cells.Item(1000, 1)
    $stage.Formula = "=TEXT(" + $digits + ",""000000000"")"
    $stage.Copy($ws.Cells.Item($row, $col))
    $stage.Clear()
}

# --- 1. Delete rows that are removed entirely in the new export ---
# (work from the bottom of the sheet upward so earlier row numbers stay valid)
$ws.Rows("36:36").Delete()   # 005000460 MARIANA   968.12
$ws.Rows("35:35").Delete()   # 005341184 BRENO      1102.66
$ws.Rows("34:34").Delete()   # 008002502 JORGEANA   1500
$ws.Rows("29:29").Delete()   # 004413537 CLAUDIA    5487.09
$ws.Rows("28:28").Delete()   # 004212476 MARIA      6420.22
$ws.Rows("27:27").Delete()   # 005046919 MARIANA    8602.79
$ws.Rows("25:25").Delete()   # 005332720 LAURA      9397.94
$ws.Rows("6:6").Delete()     # 001761119 BLUEMETRIX 202868.56 (old slot)
$ws.Rows("5:5").Delete()     # 004450724 ASSAKO     244209.73 (old slot)
$ws.Rows("4:4").Delete()     # 004452912 BRUNO      363769.47 (old slot/value)
$ws.Rows("3:3").Delete()     # 004352384 BRASFORT   422904.69

# --- 2. Update values that changed in place ---
# CAMILA is now row 2 after the deletions above
$ws.Cells.Item(2, 3).Value = 171475.78

# 004267119 / ANA is now row 13 after the deletions above
$ws.Cells.Item(13, 3).Value = 34000

# --- 3. Insert the 3 new top rows (BRUNO, ASSAKO, BLUEMETRIX) above CAMILA ---
$ws.Rows("2:4").Insert()

Set-TextAccount 2 1 4452912
$ws.Cells.Item(2, 2).Value = "BRUNO"
$ws.Cells.Item(2, 3).Value = 250000

Set-TextAccount 3 1 4450724
$ws.Cells.Item(3, 2).Value = "ASSAKO"
$ws.Cells.Item(3, 3).Value = 244209.73

Set-TextAccount 4 1 1761119
$ws.Cells.Item(4, 2).Value = "BLUEMETRIX"
$ws.Cells.Item(4, 3).Value = 204937.36
